# "Edit layout of the map"
# Moves six map-label textboxes to new positions on slide 1, and
# splits the "SAY CHEESY PRIMES" label into two lines: "DA VINCI" / "PRIMES".
#
# Note on coordinates: Shape.Left / .Top are expressed in points (1 pt =
# 12700 EMU) in the PowerPoint object model, while the OOXML stores EMUs.
# The literal values below are chosen so that the point -> EMU round trip
# lands exactly on the target EMU values from the source diff (the object
# model's internal float precision otherwise truncates by up to 1 EMU).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "NOR ORACLE": (6576110,4212100) -> (7670142,5402619)
$shp = $s.Shapes.Item("TextBox 61")
$shp.Left = 603.948218976378
$shp.Top  = 425.4030708661417

# "FACT-SEARCH": (5515571,5362277) -> (6567794,1087745)
$shp = $s.Shapes.Item("TextBox 62")
$shp.Left = 517.1491438582677
$shp.Top  = 85.6492125984252

# "MANDELE EFFECT": (6487930,1103698) -> (7548688,2168538)
$shp = $s.Shapes.Item("TextBox 63")
$shp.Left = 594.3848818897637
$shp.Top  = 170.75103362204726

# "HACK & HECK": (7531117,2198775) -> (8583487,1050428)
$shp = $s.Shapes.Item("TextBox 70")
$shp.Left = 675.8651481102363
$shp.Top  = 82.71086614173228

# "BOMBARDAIR": (8448306,1382745) -> (5285558,5468186)
$shp = $s.Shapes.Item("TextBox 71")
$shp.Left = 416.1856892913386
$shp.Top  = 430.5658467716535

# "SAY CHEESY PRIMES": (7525175,5370599) -> (6398171,4282407)
# and text "SAY CHEESY PRIMES" -> "DA VINCI" / "PRIMES" (two paragraphs)
$shp = $s.Shapes.Item("TextBox 85")
$shp.Left = 503.79299212598426
$shp.Top  = 337.19741157480314
$shp.TextFrame.TextRange.Text = "DA VINCI" + [char]13 + "PRIMES"
